# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# handback for "5c468bd5-6e1f-4ccf-bb66-5f6027a6756f.md" failed:
#   - Status changes from "Ready for handoff" to "Handback transform failed"
#     on every sheet/table row referencing that file.
#   - The "Error Detail" column (P) on the zh-cn and de-de tables is filled
#     in with a message explaining the handback/handoff file name mismatch.
#   - The "Error Detail" column is widened so the message is readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet/table ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = "Handback file name: uuqwuewd.v1f is different with handoff file name: 5c468bd5-6e1f-4ccf-bb66-5f6027a6756f.30ff0bb675f6dc0c97076405c2680b138e84b555.zh-cn."
# width="40" in the XML corresponds to a ColumnWidth of 40 minus the
# fixed 5/6 character padding that this engine adds when writing the file.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet/table ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = "Handback file name: uuqwuewd.v1f is different with handoff file name: 5c468bd5-6e1f-4ccf-bb66-5f6027a6756f.30ff0bb675f6dc0c97076405c2680b138e84b555.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
